$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$newName = "1062-MS-EPP-DB-DL-REC-NON-RNI-CTPD-DL-MD-TR-1-LateRepayment"

# Update the product name value on both sheets (space removed from "Late Repayment")
$ws1.Range("B1").Value = $newName
$ws2.Range("B1").Value = $newName

# Make ProductLoanInput the active sheet, with B1 selected on both sheets
[void]$ws1.Activate()
[void]$ws1.Range("B1").Select()
[void]$ws2.Range("B1").Select()
[void]$ws1.Activate()
